# Add a new worksheet ("1") at the end of the workbook with the
# full-obs-pendulum summary row, mirroring the existing "17" sheet layout.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "1"

$newSheet.Range("A1").Value = "ind"
$newSheet.Range("B1").Value = "fitness"
$newSheet.Range("A2").Value = "conditional(conditional(y, y), add(vel, x))"
$newSheet.Range("B2").Value = -317
